$d = $word.ActiveDocument

# --- Change 1: insert the new Canaza/SQLite bibliography reference + its Consulta paragraph
# before the existing "[Mendez, 2013]" paragraph ---
$idx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Contains("Méndez, 2013")) {
        $idx = $i
    }
}
$prevPara = $d.Paragraphs.Item($idx - 1)
$insertPos = $prevPara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" /><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /><w:szCs w:val="24" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t>[</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /><w:bCs /><w:iCs /><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>Canaza</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t>, 2014</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t xml:space="preserve">] </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /><w:bCs /><w:iCs /><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>Canaza</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /><w:bCs /><w:iCs /><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /><w:bCs /><w:iCs /><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>E</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /><w:bCs /><w:iCs /><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /><w:bCs /><w:iCs /><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /><w:bCs /><w:iCs /><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t>SQLite en Unity 3D</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t>", 2014</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t>http://es.slideshare.net/hnesys/sqlite-in-unity3d</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t>[Consulta: 29</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t xml:space="preserve"> de </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t>Agosto</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" /></w:rPr><w:t xml:space="preserve"> del 2014]</w:t></w:r></w:p>'
$insertRange.InsertXML($xml1)

Write-Output "done change 1"
